# "fixed December Rent bug"
#
# - Adds two new sheets ("2022" and "2020"), both snapshotted from the
#   original "2021" sheet data (i.e. before the December-column fix below
#   was applied), positioned right after "2021".
# - "2022" additionally gets the January (B) amounts mirrored into the
#   November (L) column, and has the December (M) rent figures cleared.
# - "2020" is a fully "fixed" copy: December (M) rent figures cleared and
#   every SUM() formula in column N / row 7 replaced by its resting value.
# - "2021" itself gets the same December-rent fix applied in place.

$wb = $excel.ActiveWorkbook
$ws2021 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Create "2022" as a copy of the untouched "2021" sheet.
# ---------------------------------------------------------------------
$ws2021.Copy([System.Reflection.Missing]::Value, $ws2021)
$ws2022 = $wb.Worksheets.Item(2)
$ws2022.Name = "2022"

$ws2022.Range("L2").Value2 = $ws2022.Range("B2").Value2
$ws2022.Range("L5").Value2 = $ws2022.Range("B5").Value2

$ws2022.Range("M3").Value2 = 0
$ws2022.Range("M4").Value2 = 0
$ws2022.Range("M6").Value2 = 0

# ---------------------------------------------------------------------
# 2) Create "2020" as another copy of the untouched "2021" sheet, then
#    apply the full December-rent-bug fix (static zeros, no formulas).
# ---------------------------------------------------------------------
$ws2021.Copy([System.Reflection.Missing]::Value, $ws2022)
$ws2020 = $wb.Worksheets.Item(3)
$ws2020.Name = "2020"

$ws2020.Range("M3").Value2 = 0
$ws2020.Range("M4").Value2 = 0
$ws2020.Range("M6").Value2 = 0

$ws2020.Range("N2").Value2 = 0
$ws2020.Range("N3").Value2 = 0
$ws2020.Range("N4").Value2 = 0
$ws2020.Range("N5").Value2 = 0
$ws2020.Range("N6").Value2 = 0

$ws2020.Range("B7").Value2 = 0
$ws2020.Range("C7").Value2 = 0
$ws2020.Range("D7").Value2 = 0
$ws2020.Range("E7").Value2 = 0
$ws2020.Range("F7").Value2 = 0
$ws2020.Range("G7").Value2 = 0
$ws2020.Range("H7").Value2 = 0
$ws2020.Range("I7").Value2 = 0
$ws2020.Range("J7").Value2 = 0
$ws2020.Range("K7").Value2 = 0
$ws2020.Range("L7").Value2 = 0
$ws2020.Range("M7").Value2 = 0
$ws2020.Range("N7").Value2 = 0

# ---------------------------------------------------------------------
# 3) Apply the same December-rent-bug fix to the original "2021" sheet.
# ---------------------------------------------------------------------
$ws2021.Range("M3").Value2 = 0
$ws2021.Range("M4").Value2 = 0
$ws2021.Range("M6").Value2 = 0

$ws2021.Range("N2").Value2 = 0
$ws2021.Range("N3").Value2 = 0
$ws2021.Range("N4").Value2 = 0
$ws2021.Range("N5").Value2 = 0
$ws2021.Range("N6").Value2 = 0

$ws2021.Range("B7").Value2 = 0
$ws2021.Range("C7").Value2 = 0
$ws2021.Range("D7").Value2 = 0
$ws2021.Range("E7").Value2 = 0
$ws2021.Range("F7").Value2 = 0
$ws2021.Range("G7").Value2 = 0
$ws2021.Range("H7").Value2 = 0
$ws2021.Range("I7").Value2 = 0
$ws2021.Range("J7").Value2 = 0
$ws2021.Range("K7").Value2 = 0
$ws2021.Range("L7").Value2 = 0
$ws2021.Range("M7").Value2 = 0
$ws2021.Range("N7").Value2 = 0

# Restore the original active/selected sheet.
$ws2021.Select()
